$d = $word.ActiveDocument

$d.Content.Find.Execute("97×96=", $true, $false, $false, $false, $false, $true, 1, $false, "59×36=", 2) | Out-Null
$d.Content.Find.Execute("62×17=", $true, $false, $false, $false, $false, $true, 1, $false, "60×31=", 2) | Out-Null
$d.Content.Find.Execute("98×55=", $true, $false, $false, $false, $false, $true, 1, $false, "46×63=", 2) | Out-Null
$d.Content.Find.Execute("59×93=", $true, $false, $false, $false, $false, $true, 1, $false, "16×13=", 2) | Out-Null
$d.Content.Find.Execute("17×81=", $true, $false, $false, $false, $false, $true, 1, $false, "56×62=", 2) | Out-Null
$d.Content.Find.Execute("19×36=", $true, $false, $false, $false, $false, $true, 1, $false, "48×85=", 2) | Out-Null
$d.Content.Find.Execute("86×27=", $true, $false, $false, $false, $false, $true, 1, $false, "31×30=", 2) | Out-Null
$d.Content.Find.Execute("17×69=", $true, $false, $false, $false, $false, $true, 1, $false, "84×76=", 2) | Out-Null
$d.Content.Find.Execute("60×39=", $true, $false, $false, $false, $false, $true, 1, $false, "16×38=", 2) | Out-Null
$d.Content.Find.Execute("67×29=", $true, $false, $false, $false, $false, $true, 1, $false, "47×67=", 2) | Out-Null
$d.Content.Find.Execute("28×49=", $true, $false, $false, $false, $false, $true, 1, $false, "27×67=", 2) | Out-Null
$d.Content.Find.Execute("77×81=", $true, $false, $false, $false, $false, $true, 1, $false, "75×31=", 2) | Out-Null
$d.Content.Find.Execute("93×15=", $true, $false, $false, $false, $false, $true, 1, $false, "47×56=", 2) | Out-Null
$d.Content.Find.Execute("57×48=", $true, $false, $false, $false, $false, $true, 1, $false, "26×89=", 2) | Out-Null
$d.Content.Find.Execute("80×86=", $true, $false, $false, $false, $false, $true, 1, $false, "87×77=", 2) | Out-Null
$d.Content.Find.Execute("77×61=", $true, $false, $false, $false, $false, $true, 1, $false, "21×25=", 2) | Out-Null
$d.Content.Find.Execute("22×19=", $true, $false, $false, $false, $false, $true, 1, $false, "63×29=", 2) | Out-Null
$d.Content.Find.Execute("37×70=", $true, $false, $false, $false, $false, $true, 1, $false, "48×60=", 2) | Out-Null
$d.Content.Find.Execute("43×88=", $true, $false, $false, $false, $false, $true, 1, $false, "97×11=", 2) | Out-Null
$d.Content.Find.Execute("45×63=", $true, $false, $false, $false, $false, $true, 1, $false, "92×17=", 2) | Out-Null
$d.Content.Find.Execute("48×40=", $true, $false, $false, $false, $false, $true, 1, $false, "38×62=", 2) | Out-Null
$d.Content.Find.Execute("51×62=", $true, $false, $false, $false, $false, $true, 1, $false, "62×80=", 2) | Out-Null
$d.Content.Find.Execute("84×88=", $true, $false, $false, $false, $false, $true, 1, $false, "55×39=", 2) | Out-Null
$d.Content.Find.Execute("65×85=", $true, $false, $false, $false, $false, $true, 1, $false, "38×31=", 2) | Out-Null
$d.Content.Find.Execute("89×51=", $true, $false, $false, $false, $false, $true, 1, $false, "36×58=", 2) | Out-Null
